$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the distinctive font/color style (applied to old E5) by copying
# it onto the new D6 cell (the row that will hold the "Cancer Center" /
# Dana Farber entry) before any values are overwritten.
$ws.Range("E5").Copy($ws.Range("D6")) | Out-Null

# Rename the worksheet to match the new, narrower data set.
$ws.Name = "page_title"

# Write the new path/type/language/page_title table.
$ws.Range("A1").Value = "path"
$ws.Range("B1").Value = "type"
$ws.Range("C1").Value = "language"
$ws.Range("D1").Value = "page_title"
$ws.Range("A2").Value = "/about-cancer/coping/feelings"
$ws.Range("B2").Value = "Article"
$ws.Range("C2").Value = "English"
$ws.Range("D2").Value = "Feelings and Cancer"
$ws.Range("A3").Value = "/espanol/cancer/sobrellevar/sentimientos/relajarse"
$ws.Range("B3").Value = "Article"
$ws.Range("C3").Value = "Spanish"
$ws.Range("D3").Value = "Aprenda a relajarse"
$ws.Range("A4").Value = "/espanol/cancer/sobrellevar/sentimientos"
$ws.Range("B4").Value = "Article"
$ws.Range("C4").Value = "Spanish"
$ws.Range("D4").Value = "Los sentimientos y el cáncer"
$ws.Range("A5").Value = "/about-cancer/coping/feelings/relaxation"
$ws.Range("B5").Value = "Article"
$ws.Range("C5").Value = "English"
$ws.Range("D5").Value = "Learning to Relax"
$ws.Range("A6").Value = "/about-cancer/coping/feelings/relaxation/dfharvard"
$ws.Range("B6").Value = "Cancer Center"
$ws.Range("C6").Value = "English"
$ws.Range("D6").Value = "Dana Farber/Harvard Cancer Center"
$ws.Range("A7").Value = "/about-cancer/coping/feelings/relaxation/loukissas-jennifer"
$ws.Range("B7").Value = "Biography"
$ws.Range("C7").Value = "English"
$ws.Range("D7").Value = "Jennifer K. Loukissas, M.P.P."
$ws.Range("A8").Value = "/news-events/cancer-currents-blog/2019/vitamin-d-supplement-cancer-prevention"
$ws.Range("B8").Value = "Blog"
$ws.Range("C8").Value = "English"
$ws.Range("D8").Value = "Vitamin D Supplements Don’t Reduce Cancer Incidence, Trial Shows"
$ws.Range("A9").Value = "/espanol/noticias/temas-y-relatos-blog/2019/vitamina-d-complemento-cancer-prevencion"
$ws.Range("B9").Value = "Blog"
$ws.Range("C9").Value = "Spanish"
$ws.Range("D9").Value = "Estudio indica que complementos de vitamina D no reducen la incidencia de cáncer"
$ws.Range("A10").Value = "/espanol/noticias/comunicados-de-prensa/2018/leucemia-llc-ibrutinib-estudio"
$ws.Range("B10").Value = "Press Release"
$ws.Range("C10").Value = "Spanish"
$ws.Range("D10").Value = "Leucemia - Ibrutinib más rituximab fue superior al tratamiento convencional para algunos pacientes con leucemia linfocítica crónica"
$ws.Range("A11").Value = "/news-events/press-releases/2018/leukemia-cll-ibrutinib-trial"
$ws.Range("B11").Value = "Press Release"
$ws.Range("C11").Value = "English"
$ws.Range("D11").Value = "Leukemia - Ibrutinib plus rituximab superior to standard treatment for some patients with chronic leukemia"

# Drop the now-unused 5th (E) column that held the old browser_title data.
$ws.Range("E1:E11").Clear() | Out-Null

# Column widths to roughly match the refreshed layout.
$ws.Columns.Item(1).ColumnWidth = 49
$ws.Columns.Item(2).ColumnWidth = 12.5
$ws.Columns.Item(4).ColumnWidth = 12.5

# Selection / active cell below the data, matching the saved view state.
$ws.Range("D12").Select() | Out-Null
